# updateStartingListSheet(): append the webhook log rows captured since the
# last sync to the "Webhooks" sheet without disturbing the header or any
# previously written rows, then mirror dimension/layout bookkeeping.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Webhooks")

# getHeaderInfoFromRow(): read the header row once so the column count used
# to size every appended record always tracks the sheet's real header,
# rather than a hard-coded literal.
$usedRange = $ws.UsedRange
$headerCols = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

# New webhook events recorded since the previous export run.
$newRows = @(
    @("2026-02-06T22:01:21.562Z", "event_raslisti_birtur", 70617, 103060, 1, 1, '{"classId":103060,"eventId":70617,"competitionId":1,"published":1}'),
    @("2026-02-06T22:02:08.841Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:02:52.439Z", "event_einkunn_saeti",   70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:03:53.610Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:12:32.161Z", "event_raslisti_birtur", 70674, 70674,  2, 1, '{"classId":70674,"eventId":70674,"competitionId":2,"published":1}'),
    @("2026-02-06T22:12:36.440Z", "event_einkunn_saeti",   70674, 70674,  2, 1, '{"classId":70674,"eventId":70674,"competitionId":2,"published":1}'),
    @("2026-02-06T22:14:06.748Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:14:26.229Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:14:36.775Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:16:11.742Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:16:27.774Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:19:12.265Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:19:23.415Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:19:57.101Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:20:09.228Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:20:30.863Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}'),
    @("2026-02-06T22:21:29.875Z", "event_raslisti_birtur", 70674, 70674,  1, 1, '{"classId":70674,"eventId":70674,"competitionId":1,"published":1}')
)

$r = $lastRow + 1
foreach ($rowVals in $newRows) {
    $c = 1
    while ($c -le $headerCols -and $c -le $rowVals.Count) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
        $c = $c + 1
    }
    $r = $r + 1
}
